$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "iaest-measure:actividades-desarrollo-rural"

# Row 3
$ws.Range("C3").Value = "dim"
$ws.Range("H3").Value = "medida"

# Row 4
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("H4").Value = "xsd:int"

# Row 5 - remove H5 entirely
$ws.Range("H5").Clear()
